$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the date number format used by column D so the new/edited
# cells keep the same display format (style index 2 in the original file).
$dateFormat = $ws.Range("D14").NumberFormat

# --- Insert the new row 15, a copy of the original row 14 data ---
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 45015
$ws.Range("D15").NumberFormat = $dateFormat
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112041
$ws.Range("G15").Value = "Fruto del paraíso"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 24000
$ws.Range("N15").Value = "$/caja 18 kilos empedrada"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 1333
$ws.Range("Q15").Value = 18
$ws.Range("R15").Value = "Hortaliza"

# --- Update row 14 with the new weekly values ---
$ws.Range("D14").Value = 45041
$ws.Range("D14").NumberFormat = $dateFormat
$ws.Range("J14").Value = 80
